$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44559
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100101
$ws.Range("H4").Value = "Berries"
$ws.Range("I4").Value = 100101008
$ws.Range("J4").Value = "Mora"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6500
$ws.Range("Q4").Value = "$/bandeja 2 kilos"
$ws.Range("R4").Value = "Región de Ñuble"
$ws.Range("S4").Value = 3250
$ws.Range("T4").Value = 2

$ws.Range("A5").Value = 11
$ws.Range("B5").Value = "Vega Monumental Concepción"
$ws.Range("C5").Value = "Bíobío"
$ws.Range("D5").Value = 44559
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100101
$ws.Range("H5").Value = "Berries"
$ws.Range("I5").Value = 100101008
$ws.Range("J5").Value = "Mora"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 5000
$ws.Range("O5").Value = 5000
$ws.Range("P5").Value = 5000
$ws.Range("Q5").Value = "$/bandeja 2 kilos"
$ws.Range("R5").Value = "Región de Ñuble"
$ws.Range("S5").Value = 2500
$ws.Range("T5").Value = 2

# Apply the same style as D2/D3 (date format) to D4/D5
$ws.Range("D2").Copy()
$ws.Range("D4:D5").PasteSpecial(-4122)
